$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing Luca "Documentazione esterna/interna" amounts (rows 43-45) ---
$ws.Range("C43").Value = 180
$ws.Range("C44").Value = 180
$ws.Range("C45").Value = 120

# --- New row 46: Giovanni / Documentazione interna / 300 / 12/15/2018 (text date) ---
$ws.Range("A46").Value = "Giovanni"
$ws.Range("B46").Value = "Documentazione interna"
$ws.Range("C46").Value = 300
$ws.Range("D42").Copy()
$ws.Range("D46").PasteSpecial(-4122)
$ws.Range("D46").Formula = '="12/15/2018"'
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)

# --- New row 47: Giovanni / Documentazione esterna / 150 / 12/16/2018 (text date) ---
$ws.Range("A47").Value = "Giovanni"
$ws.Range("B47").Value = "Documentazione esterna"
$ws.Range("C47").Value = 150
$ws.Range("D42").Copy()
$ws.Range("D47").PasteSpecial(-4122)
$ws.Range("D47").Formula = '="12/16/2018"'
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)

# --- New row 48: Luca / Documentazione interna / 82 / 03/04/2019 (real date) ---
$ws.Range("A48").Value = "Luca"
$ws.Range("B48").Value = "Documentazione interna"
$ws.Range("C48").Value = 82
$ws.Range("D48").Value = 43528
$ws.Range("D42").Copy()
$ws.Range("D48").PasteSpecial(-4122)

# --- New row 49: Hristina / Documentazione interna / 82 / 03/04/2019 (real date) ---
$ws.Range("A49").Value = "Hristina"
$ws.Range("B49").Value = "Documentazione interna"
$ws.Range("C49").Value = 82
$ws.Range("D49").Value = 43528
$ws.Range("D42").Copy()
$ws.Range("D49").PasteSpecial(-4122)

# --- New row 50: Giovanni / Documentazione interna / 82 / 03/04/2019 (real date) ---
$ws.Range("A50").Value = "Giovanni"
$ws.Range("B50").Value = "Documentazione interna"
$ws.Range("C50").Value = 82
$ws.Range("D50").Value = 43528
$ws.Range("D42").Copy()
$ws.Range("D50").PasteSpecial(-4122)

# --- Update selection to mirror the author's last position ---
$ws.Range("C50").Select()
